$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F ("dSF") values to match repulled data / mean calculation fix
$ws.Range("F2").Value = -9
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -8
$ws.Range("F13").Value = -5
$ws.Range("F16").Value = -2
